# Auto-generated market-data refresh edit.
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ (columns H-N) across all
# 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect refreshed market
# board data. Some rows gain or lose an HQ-profit (N) cell when HQ pricing
# becomes available/unavailable for that item.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 72.0625
$ws.Cells.Item(2, 9).Value = 67.69231000000001
$ws.Cells.Item(2, 11).Value = 67.69231000000001
$ws.Cells.Item(2, 13).Value = 45.30768999999999
$ws.Cells.Item(12, 8).Value = 366.66666
$ws.Cells.Item(12, 9).Value = 366.66666
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 366.66666
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 13).Value = -196.66666
$ws.Cells.Item(12, 14).ClearContents()
$ws.Cells.Item(17, 8).Value = 361.41177
$ws.Cells.Item(17, 10).Value = 368.21875
$ws.Cells.Item(17, 12).Value = 1104.65625
$ws.Cells.Item(17, 14).Value = -1440.65625
$ws.Cells.Item(32, 8).Value = 1131.5
$ws.Cells.Item(32, 10).Value = 1078.8572
$ws.Cells.Item(32, 12).Value = 1078.8572
$ws.Cells.Item(32, 14).Value = -1730.8572
$ws.Cells.Item(40, 8).Value = 1969863.1
$ws.Cells.Item(40, 9).Value = 13504.444
$ws.Cells.Item(40, 10).Value = 4170766.8
$ws.Cells.Item(40, 11).Value = 13504.444
$ws.Cells.Item(40, 12).Value = 4170766.8
$ws.Cells.Item(40, 13).Value = -13329.444
$ws.Cells.Item(40, 14).Value = -4171116.8
$ws.Cells.Item(51, 8).Value = 10599
$ws.Cells.Item(51, 9).Value = 5000
$ws.Cells.Item(51, 11).Value = 5000
$ws.Cells.Item(51, 13).Value = -4516
$ws.Cells.Item(74, 8).Value = 35723332
$ws.Cells.Item(74, 9).Value = 125002936
$ws.Cells.Item(74, 11).Value = 125002936
$ws.Cells.Item(74, 13).Value = -125002000
$ws.Cells.Item(77, 8).Value = 35723332
$ws.Cells.Item(77, 9).Value = 125002936
$ws.Cells.Item(77, 11).Value = 625014680
$ws.Cells.Item(77, 13).Value = -625010000
$ws.Cells.Item(106, 8).Value = 3832.875
$ws.Cells.Item(106, 9).Value = 3832.875
$ws.Cells.Item(106, 11).Value = 3832.875
$ws.Cells.Item(106, 13).Value = -3201.875
$ws.Cells.Item(132, 8).Value = 2515.7917
$ws.Cells.Item(132, 9).Value = 2515.7917
$ws.Cells.Item(132, 11).Value = 7547.375100000001
$ws.Cells.Item(132, 13).Value = -5017.375100000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1814.0714
$ws.Cells.Item(2, 9).Value = 1101.1818
$ws.Cells.Item(2, 11).Value = 1101.1818
$ws.Cells.Item(2, 13).Value = -988.1818000000001
$ws.Cells.Item(32, 8).Value = 2157431.8
$ws.Cells.Item(32, 9).Value = 2234378
$ws.Cells.Item(32, 11).Value = 2234378
$ws.Cells.Item(32, 13).Value = -2234091
$ws.Cells.Item(74, 8).Value = 44293.043
$ws.Cells.Item(74, 9).Value = 60061.35
$ws.Cells.Item(74, 10).Value = 5998.5713
$ws.Cells.Item(74, 11).Value = 60061.35
$ws.Cells.Item(74, 12).Value = 5998.5713
$ws.Cells.Item(74, 13).Value = -59187.35
$ws.Cells.Item(74, 14).Value = -7746.5713
$ws.Cells.Item(77, 8).Value = 44293.043
$ws.Cells.Item(77, 9).Value = 60061.35
$ws.Cells.Item(77, 10).Value = 5998.5713
$ws.Cells.Item(77, 11).Value = 300306.75
$ws.Cells.Item(77, 12).Value = 29992.8565
$ws.Cells.Item(77, 13).Value = -295938.75
$ws.Cells.Item(77, 14).Value = -38728.85649999999
$ws.Cells.Item(97, 8).Value = 3974441.5
$ws.Cells.Item(97, 9).Value = 457.41666
$ws.Cells.Item(97, 11).Value = 457.41666
$ws.Cells.Item(97, 13).Value = 38.58334000000002
$ws.Cells.Item(116, 8).Value = 1814.0714
$ws.Cells.Item(116, 9).Value = 1101.1818
$ws.Cells.Item(116, 11).Value = 1101.1818
$ws.Cells.Item(116, 13).Value = 1192.8182
$ws.Cells.Item(122, 8).Value = 4681.8667
$ws.Cells.Item(122, 10).Value = 6318.1665
$ws.Cells.Item(122, 12).Value = 18954.4995
$ws.Cells.Item(122, 14).Value = -23854.4995
$ws.Cells.Item(132, 8).Value = 9927.786
$ws.Cells.Item(132, 9).Value = 11227
$ws.Cells.Item(132, 10).Value = 9087.117
$ws.Cells.Item(132, 11).Value = 33681
$ws.Cells.Item(132, 12).Value = 27261.351
$ws.Cells.Item(132, 13).Value = -31151
$ws.Cells.Item(132, 14).Value = -32321.351

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1814.0714
$ws.Cells.Item(3, 9).Value = 1101.1818
$ws.Cells.Item(3, 11).Value = 1101.1818
$ws.Cells.Item(3, 13).Value = -987.1818000000001
$ws.Cells.Item(82, 8).Value = 0
$ws.Cells.Item(82, 9).Value = 0
$ws.Cells.Item(82, 11).Value = 0
$ws.Cells.Item(82, 13).ClearContents()
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 13).ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(5, 8).Value = 404
$ws.Cells.Item(5, 9).Value = 404
$ws.Cells.Item(5, 11).Value = 404
$ws.Cells.Item(5, 13).Value = -292
$ws.Cells.Item(31, 8).Value = 7920.46
$ws.Cells.Item(31, 9).Value = 3867.7144
$ws.Cells.Item(31, 10).Value = 9496.527
$ws.Cells.Item(31, 11).Value = 3867.7144
$ws.Cells.Item(31, 12).Value = 9496.527
$ws.Cells.Item(31, 13).Value = -3572.7144
$ws.Cells.Item(31, 14).Value = -10086.527
$ws.Cells.Item(34, 8).Value = 7920.46
$ws.Cells.Item(34, 9).Value = 3867.7144
$ws.Cells.Item(34, 10).Value = 9496.527
$ws.Cells.Item(34, 11).Value = 3867.7144
$ws.Cells.Item(34, 12).Value = 9496.527
$ws.Cells.Item(34, 13).Value = -3665.7144
$ws.Cells.Item(34, 14).Value = -9900.527
$ws.Cells.Item(107, 8).Value = 1461.5264
$ws.Cells.Item(107, 9).Value = 1127.1072
$ws.Cells.Item(107, 10).Value = 2397.9
$ws.Cells.Item(107, 11).Value = 1127.1072
$ws.Cells.Item(107, 12).Value = 2397.9
$ws.Cells.Item(107, 13).Value = 792.8928000000001
$ws.Cells.Item(107, 14).Value = -6237.9

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 4447066.5
$ws.Cells.Item(5, 9).Value = 8000920.5
$ws.Cells.Item(5, 10).Value = 4749
$ws.Cells.Item(5, 11).Value = 24002761.5
$ws.Cells.Item(5, 12).Value = 14247
$ws.Cells.Item(5, 13).Value = -24002649.5
$ws.Cells.Item(5, 14).Value = -14471
$ws.Cells.Item(12, 8).Value = 656.4
$ws.Cells.Item(12, 9).Value = 676.7143
$ws.Cells.Item(12, 10).Value = 645.46155
$ws.Cells.Item(12, 11).Value = 2030.1429
$ws.Cells.Item(12, 12).Value = 1936.38465
$ws.Cells.Item(12, 13).Value = -1857.1429
$ws.Cells.Item(12, 14).Value = -2282.38465
$ws.Cells.Item(122, 8).Value = 1490651.5
$ws.Cells.Item(122, 9).Value = 3537497.5
$ws.Cells.Item(122, 10).Value = 2036.091
$ws.Cells.Item(122, 11).Value = 31837477.5
$ws.Cells.Item(122, 12).Value = 18324.819
$ws.Cells.Item(122, 13).Value = -31835027.5
$ws.Cells.Item(122, 14).Value = -23224.819
$ws.Cells.Item(131, 8).Value = 1614.8096
$ws.Cells.Item(131, 9).Value = 830
$ws.Cells.Item(131, 10).Value = 2661.2222
$ws.Cells.Item(131, 11).Value = 2490
$ws.Cells.Item(131, 12).Value = 7983.6666
$ws.Cells.Item(131, 13).Value = 2550
$ws.Cells.Item(131, 14).Value = -18063.6666
$ws.Cells.Item(132, 8).Value = 8070.5454
$ws.Cells.Item(132, 9).Value = 6984.5
$ws.Cells.Item(132, 10).Value = 8975.583000000001
$ws.Cells.Item(132, 11).Value = 62860.5
$ws.Cells.Item(132, 12).Value = 80780.247
$ws.Cells.Item(132, 13).Value = -60330.5
$ws.Cells.Item(132, 14).Value = -85840.247
$ws.Cells.Item(135, 8).Value = 4447066.5
$ws.Cells.Item(135, 9).Value = 8000920.5
$ws.Cells.Item(135, 10).Value = 4749
$ws.Cells.Item(135, 11).Value = 72008284.5
$ws.Cells.Item(135, 12).Value = 42741
$ws.Cells.Item(135, 13).Value = -72005749.5
$ws.Cells.Item(135, 14).Value = -47811

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 1111184.4
$ws.Cells.Item(2, 9).Value = 65.90909000000001
$ws.Cells.Item(2, 10).Value = 2857227.8
$ws.Cells.Item(2, 11).Value = 65.90909000000001
$ws.Cells.Item(2, 12).Value = 2857227.8
$ws.Cells.Item(2, 13).Value = 47.09090999999999
$ws.Cells.Item(2, 14).Value = -2857453.8
$ws.Cells.Item(15, 8).Value = 56000
$ws.Cells.Item(15, 10).Value = 56000
$ws.Cells.Item(15, 12).Value = 56000
$ws.Cells.Item(15, 14).Value = -56576
$ws.Cells.Item(80, 8).Value = 35104.965
$ws.Cells.Item(80, 9).Value = 1678.591
$ws.Cells.Item(80, 11).Value = 1678.591
$ws.Cells.Item(80, 13).Value = -680.5909999999999
$ws.Cells.Item(81, 8).Value = 56000
$ws.Cells.Item(81, 10).Value = 56000
$ws.Cells.Item(81, 12).Value = 56000
$ws.Cells.Item(81, 14).Value = -57996
$ws.Cells.Item(82, 8).Value = 0
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 14).ClearContents()
$ws.Cells.Item(83, 8).Value = 35104.965
$ws.Cells.Item(83, 9).Value = 1678.591
$ws.Cells.Item(83, 11).Value = 8392.955
$ws.Cells.Item(83, 13).Value = -3400.955
$ws.Cells.Item(84, 8).Value = 56000
$ws.Cells.Item(84, 10).Value = 56000
$ws.Cells.Item(84, 12).Value = 168000
$ws.Cells.Item(84, 14).Value = -177984
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 7857.484
$ws.Cells.Item(113, 9).Value = 5121.769
$ws.Cells.Item(113, 11).Value = 5121.769
$ws.Cells.Item(113, 13).Value = -2951.769
$ws.Cells.Item(132, 8).Value = 4199.1934
$ws.Cells.Item(132, 9).Value = 1861.35
$ws.Cells.Item(132, 10).Value = 8449.817999999999
$ws.Cells.Item(132, 11).Value = 5584.049999999999
$ws.Cells.Item(132, 12).Value = 25349.454
$ws.Cells.Item(132, 13).Value = -3054.049999999999
$ws.Cells.Item(132, 14).Value = -30409.454

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2566.1875
$ws.Cells.Item(22, 10).Value = 3531.25
$ws.Cells.Item(22, 12).Value = 3531.25
$ws.Cells.Item(22, 14).Value = -4121.25
$ws.Cells.Item(27, 8).Value = 2566.1875
$ws.Cells.Item(27, 10).Value = 3531.25
$ws.Cells.Item(27, 12).Value = 3531.25
$ws.Cells.Item(27, 14).Value = -3745.25
$ws.Cells.Item(46, 8).Value = 1380842.9
$ws.Cells.Item(46, 9).Value = 3135563.2
$ws.Cells.Item(46, 10).Value = 2134
$ws.Cells.Item(46, 11).Value = 3135563.2
$ws.Cells.Item(46, 12).Value = 2134
$ws.Cells.Item(46, 13).Value = -3135375.2
$ws.Cells.Item(46, 14).Value = -2510

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 20007050
$ws.Cells.Item(81, 9).Value = 1312.5
$ws.Cells.Item(81, 11).Value = 2625
$ws.Cells.Item(81, 13).Value = -1564
$ws.Cells.Item(84, 8).Value = 20007050
$ws.Cells.Item(84, 9).Value = 1312.5
$ws.Cells.Item(84, 11).Value = 13125
$ws.Cells.Item(84, 13).Value = -7821
$ws.Cells.Item(113, 8).Value = 12726.954
$ws.Cells.Item(113, 9).Value = 28748.555
$ws.Cells.Item(113, 10).Value = 1635.0769
$ws.Cells.Item(113, 11).Value = 86245.66500000001
$ws.Cells.Item(113, 12).Value = 4905.2307
$ws.Cells.Item(113, 13).Value = -84075.66500000001
$ws.Cells.Item(113, 14).Value = -9245.2307
$ws.Cells.Item(136, 8).Value = 50054376
$ws.Cells.Item(136, 9).Value = 62502656
$ws.Cells.Item(136, 11).Value = 187507968
$ws.Cells.Item(136, 13).Value = -187505418
